$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH209"
$ws.Range("C2").Value = "TEXTE DIENSTE AS VEBERSEE"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24D | GRAP COUNT NUMER: NONE"
